$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = 'Última actualización: 16:34:19'
$ws1.Cells.Item(3,1).Value = 'Total filas: 328'
$ws1.Cells.Item(137,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(138,3).Value = '15_ABASTO'
$ws1.Cells.Item(178,1).Value = '11:46:46'
$ws1.Cells.Item(178,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(178,4).Value = 48
$ws1.Cells.Item(179,1).Value = '11:17:39'
$ws1.Cells.Item(179,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(179,4).Value = 77
$ws1.Cells.Item(186,1).Value = '11:46:46'
$ws1.Cells.Item(186,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(186,4).Value = 55
$ws1.Cells.Item(187,1).Value = '11:17:39'
$ws1.Cells.Item(187,3).Value = '10_OLMOS'
$ws1.Cells.Item(187,4).Value = 84
$ws1.Cells.Item(198,1).Value = '12:01:11'
$ws1.Cells.Item(198,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(198,4).Value = 66
$ws1.Cells.Item(199,1).Value = '12:50:41'
$ws1.Cells.Item(199,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(199,4).Value = 17
$ws1.Cells.Item(220,1).Value = '13:51:32'
$ws1.Cells.Item(220,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(220,4).Value = 0
$ws1.Cells.Item(221,1).Value = '12:01:11'
$ws1.Cells.Item(221,3).Value = '215A_EL PATO'
$ws1.Cells.Item(221,4).Value = 110
$ws1.Cells.Item(261,1).Value = '15:36:13'
$ws1.Cells.Item(261,3).Value = '10_OLMOS'
$ws1.Cells.Item(261,4).Value = 0
$ws1.Cells.Item(262,1).Value = '14:20:49'
$ws1.Cells.Item(262,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(262,4).Value = 76
$ws1.Cells.Item(295,1).Value = '16:34:19'
$ws1.Cells.Item(295,2).Value = '16:34'
$ws1.Cells.Item(295,4).Value = 0
$ws1.Cells.Item(296,1).Value = '16:34:19'
$ws1.Cells.Item(296,2).Value = '16:34'
$ws1.Cells.Item(296,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(296,4).Value = 0
$ws1.Cells.Item(297,1).Value = '15:59:02'
$ws1.Cells.Item(297,2).Value = '16:35'
$ws1.Cells.Item(297,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(297,4).Value = 36
$ws1.Cells.Item(298,1).Value = '16:20:15'
$ws1.Cells.Item(298,2).Value = '16:35'
$ws1.Cells.Item(298,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(298,4).Value = 15
$ws1.Cells.Item(299,1).Value = '15:36:13'
$ws1.Cells.Item(299,2).Value = '16:37'
$ws1.Cells.Item(299,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(299,4).Value = 61
$ws1.Cells.Item(300,1).Value = '15:36:13'
$ws1.Cells.Item(300,2).Value = '16:40'
$ws1.Cells.Item(300,3).Value = '17_ROMERO'
$ws1.Cells.Item(300,4).Value = 64
$ws1.Cells.Item(301,1).Value = '14:59:23'
$ws1.Cells.Item(301,2).Value = '16:42'
$ws1.Cells.Item(301,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(301,4).Value = 103
$ws1.Cells.Item(302,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(303,2).Value = '16:43'
$ws1.Cells.Item(303,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(303,4).Value = 67
$ws1.Cells.Item(304,1).Value = '14:49:07'
$ws1.Cells.Item(304,2).Value = '16:43'
$ws1.Cells.Item(304,3).Value = '225_GOMEZ'
$ws1.Cells.Item(304,4).Value = 114
$ws1.Cells.Item(305,1).Value = '15:36:13'
$ws1.Cells.Item(305,2).Value = '16:48'
$ws1.Cells.Item(305,3).Value = '15_ABASTO'
$ws1.Cells.Item(305,4).Value = 72
$ws1.Cells.Item(306,1).Value = '16:34:19'
$ws1.Cells.Item(306,2).Value = '16:50'
$ws1.Cells.Item(306,3).Value = '14_ABASTO'
$ws1.Cells.Item(306,4).Value = 16
$ws1.Cells.Item(307,1).Value = '15:59:02'
$ws1.Cells.Item(307,2).Value = '16:51'
$ws1.Cells.Item(307,3).Value = '14_ABASTO'
$ws1.Cells.Item(307,4).Value = 52
$ws1.Cells.Item(308,1).Value = '14:59:23'
$ws1.Cells.Item(308,2).Value = '16:56'
$ws1.Cells.Item(308,3).Value = '17_179 Y 38'
$ws1.Cells.Item(308,4).Value = 117
$ws1.Cells.Item(309,2).Value = '16:57'
$ws1.Cells.Item(309,3).Value = '10_OLMOS'
$ws1.Cells.Item(309,4).Value = 58
$ws1.Cells.Item(310,1).Value = '16:34:19'
$ws1.Cells.Item(310,2).Value = '17:04'
$ws1.Cells.Item(310,3).Value = '215A_EL PATO'
$ws1.Cells.Item(310,4).Value = 30
$ws1.Cells.Item(311,1).Value = '15:36:13'
$ws1.Cells.Item(311,2).Value = '17:05'
$ws1.Cells.Item(311,3).Value = '215A_EL PATO'
$ws1.Cells.Item(311,4).Value = 89
$ws1.Cells.Item(312,1).Value = '16:20:15'
$ws1.Cells.Item(312,2).Value = '17:05'
$ws1.Cells.Item(312,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(312,4).Value = 45
$ws1.Cells.Item(313,1).Value = '16:34:19'
$ws1.Cells.Item(313,2).Value = '17:10'
$ws1.Cells.Item(313,3).Value = '10_OLMOS'
$ws1.Cells.Item(313,4).Value = 36
$ws1.Cells.Item(314,1).Value = '16:34:19'
$ws1.Cells.Item(314,2).Value = '17:16'
$ws1.Cells.Item(314,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(314,4).Value = 42
$ws1.Cells.Item(315,2).Value = '17:17'
$ws1.Cells.Item(315,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(315,4).Value = 78
$ws1.Cells.Item(316,1).Value = '15:36:13'
$ws1.Cells.Item(316,2).Value = '17:21'
$ws1.Cells.Item(316,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(316,4).Value = 105
$ws1.Cells.Item(317,1).Value = '16:20:15'
$ws1.Cells.Item(317,2).Value = '17:21'
$ws1.Cells.Item(317,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(317,4).Value = 61
$ws1.Cells.Item(318,1).Value = '15:36:13'
$ws1.Cells.Item(318,2).Value = '17:24'
$ws1.Cells.Item(318,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(318,4).Value = 108
$ws1.Cells.Item(319,1).Value = '16:34:19'
$ws1.Cells.Item(319,2).Value = '17:28'
$ws1.Cells.Item(319,3).Value = '14_ABASTO'
$ws1.Cells.Item(319,4).Value = 54
$ws1.Cells.Item(320,1).Value = '16:34:19'
$ws1.Cells.Item(320,2).Value = '17:31'
$ws1.Cells.Item(320,3).Value = '15_ABASTO'
$ws1.Cells.Item(320,4).Value = 57
$ws1.Cells.Item(320,5).Value = 'LP1912'
$ws1.Cells.Item(321,1).Value = '16:20:15'
$ws1.Cells.Item(321,2).Value = '17:36'
$ws1.Cells.Item(321,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(321,4).Value = 76
$ws1.Cells.Item(321,5).Value = 'LP1912'
$ws1.Cells.Item(322,1).Value = '15:59:02'
$ws1.Cells.Item(322,2).Value = '17:37'
$ws1.Cells.Item(322,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(322,4).Value = 98
$ws1.Cells.Item(322,5).Value = 'LP1912'
$ws1.Cells.Item(323,1).Value = '15:59:02'
$ws1.Cells.Item(323,2).Value = '17:38'
$ws1.Cells.Item(323,3).Value = '17_ROMERO'
$ws1.Cells.Item(323,4).Value = 99
$ws1.Cells.Item(323,5).Value = 'LP1912'
$ws1.Cells.Item(324,1).Value = '16:34:19'
$ws1.Cells.Item(324,2).Value = '17:39'
$ws1.Cells.Item(324,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(324,4).Value = 65
$ws1.Cells.Item(324,5).Value = 'LP1912'
$ws1.Cells.Item(325,1).Value = '15:59:02'
$ws1.Cells.Item(325,2).Value = '17:40'
$ws1.Cells.Item(325,3).Value = '215B_EL PATO'
$ws1.Cells.Item(325,4).Value = 101
$ws1.Cells.Item(325,5).Value = 'LP1912'
$ws1.Cells.Item(326,1).Value = '16:34:19'
$ws1.Cells.Item(326,2).Value = '17:41'
$ws1.Cells.Item(326,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(326,4).Value = 67
$ws1.Cells.Item(326,5).Value = 'LP1912'
$ws1.Cells.Item(327,1).Value = '16:34:19'
$ws1.Cells.Item(327,2).Value = '17:50'
$ws1.Cells.Item(327,3).Value = '16_P MOR-167 Y 521'
$ws1.Cells.Item(327,4).Value = 76
$ws1.Cells.Item(327,5).Value = 'LP1912'
$ws1.Cells.Item(328,1).Value = '15:59:02'
$ws1.Cells.Item(328,2).Value = '17:51'
$ws1.Cells.Item(328,3).Value = '16_P MOR-167 Y 521'
$ws1.Cells.Item(328,4).Value = 112
$ws1.Cells.Item(328,5).Value = 'LP1912'
$ws1.Cells.Item(329,1).Value = '15:59:02'
$ws1.Cells.Item(329,2).Value = '17:52'
$ws1.Cells.Item(329,3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(329,4).Value = 113
$ws1.Cells.Item(329,5).Value = 'LP1912'
$ws1.Cells.Item(330,1).Value = '16:20:15'
$ws1.Cells.Item(330,2).Value = '18:04'
$ws1.Cells.Item(330,3).Value = '17_ROMERO'
$ws1.Cells.Item(330,4).Value = 104
$ws1.Cells.Item(330,5).Value = 'LP1912'
$ws1.Cells.Item(331,1).Value = '16:34:19'
$ws1.Cells.Item(331,2).Value = '18:21'
$ws1.Cells.Item(331,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(331,4).Value = 107
$ws1.Cells.Item(331,5).Value = 'LP1912'
$ws1.Cells.Item(332,1).Value = '16:34:19'
$ws1.Cells.Item(332,2).Value = '18:28'
$ws1.Cells.Item(332,3).Value = '215C_EL PATO'
$ws1.Cells.Item(332,4).Value = 114
$ws1.Cells.Item(332,5).Value = 'LP1912'
$ws1.Cells.Item(333,1).Value = '16:34:19'
$ws1.Cells.Item(333,2).Value = '18:32'
$ws1.Cells.Item(333,3).Value = '11X44_ETCHEVERRY'
$ws1.Cells.Item(333,4).Value = 118
$ws1.Cells.Item(333,5).Value = 'LP1912'

# --- Sheet: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = 'Última actualización: 16:34:19'
$ws2.Cells.Item(3,1).Value = 'Total filas: 37'
$ws2.Cells.Item(39,1).Value = '16:34:19'
$ws2.Cells.Item(39,2).Value = '17:04'
$ws2.Cells.Item(39,4).Value = 30
$ws2.Cells.Item(40,1).Value = '15:36:13'
$ws2.Cells.Item(40,2).Value = '17:05'
$ws2.Cells.Item(40,3).Value = '215A_EL PATO'
$ws2.Cells.Item(40,4).Value = 89
$ws2.Cells.Item(41,1).Value = '15:59:02'
$ws2.Cells.Item(41,2).Value = '17:40'
$ws2.Cells.Item(41,3).Value = '215B_EL PATO'
$ws2.Cells.Item(41,4).Value = 101
$ws2.Cells.Item(41,5).Value = 'LP1912'
$ws2.Cells.Item(42,1).Value = '16:34:19'
$ws2.Cells.Item(42,2).Value = '18:28'
$ws2.Cells.Item(42,3).Value = '215C_EL PATO'
$ws2.Cells.Item(42,4).Value = 114
$ws2.Cells.Item(42,5).Value = 'LP1912'

# --- Sheet: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = 'Última actualización: 16:34:19'
